# Auto-generated edit script: updates Leve market-price / profit
# columns (H-N) on each profession sheet to match refreshed market data.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 140.42857
$ws.Range("I5").Value = 143.83333
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 143.83333
$ws.Range("L5").Value = 120
$ws.Range("M5").Value = -28.83332999999999
$ws.Range("N5").Value = -350
# Row 6
$ws.Range("H6").Value = 180.72728
$ws.Range("I6").Value = 192.25
$ws.Range("K6").Value = 576.75
$ws.Range("M6").Value = -464.75
# Row 20
$ws.Range("H20").Value = 3840.3333
$ws.Range("J20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10460
# Row 33
$ws.Range("H33").Value = 319.6
$ws.Range("I33").Value = 301
$ws.Range("K33").Value = 301
$ws.Range("M33").Value = -72
# Row 35
$ws.Range("H35").Value = 3840.3333
$ws.Range("J35").Value = 10000
$ws.Range("L35").Value = 10000
$ws.Range("N35").Value = -10758
# Row 76
$ws.Range("H76").Value = 2000
$ws.Range("I76").Value = 2000
$ws.Range("K76").Value = 2000
$ws.Range("M76").Value = -1685
# Row 79
$ws.Range("H79").Value = 2000
$ws.Range("I79").Value = 2000
$ws.Range("K79").Value = 2000
$ws.Range("M79").Value = -908

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 121
$ws.Range("I4").Value = 105.8
$ws.Range("K4").Value = 105.8
$ws.Range("M4").Value = 10.2
# Row 5
$ws.Range("H5").Value = 98.333336
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
# Row 45
$ws.Range("H45").Value = 3600
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 3600
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 3600
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -4354
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()
# Row 74
$ws.Range("H74").Value = 791.38464
$ws.Range("I74").Value = 572.63635
$ws.Range("K74").Value = 572.63635
$ws.Range("M74").Value = 301.36365
# Row 77
$ws.Range("H77").Value = 791.38464
$ws.Range("I77").Value = 572.63635
$ws.Range("K77").Value = 2863.18175
$ws.Range("M77").Value = 1504.81825
# Row 136
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 98.333336
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
# Row 35
$ws.Range("H35").Value = 19399
$ws.Range("I35").Value = 20000
$ws.Range("J35").Value = 19158.6
$ws.Range("K35").Value = 20000
$ws.Range("L35").Value = 19158.6
$ws.Range("M35").Value = -19690
$ws.Range("N35").Value = -19778.6
# Row 80
$ws.Range("H80").Value = 299.1111
$ws.Range("J80").Value = 362.6
$ws.Range("L80").Value = 362.6
$ws.Range("N80").Value = -2358.6
# Row 83
$ws.Range("H83").Value = 299.1111
$ws.Range("J83").Value = 362.6
$ws.Range("L83").Value = 1813
$ws.Range("N83").Value = -11797
# Row 99
$ws.Range("H99").Value = 2096.375
$ws.Range("I99").Value = 2110.2856
$ws.Range("K99").Value = 2110.2856
$ws.Range("M99").Value = -612.2856000000002
# Row 134
$ws.Range("H134").Value = 5999.8
$ws.Range("I134").Value = 4999.6665
$ws.Range("J134").Value = 7500
$ws.Range("K134").Value = 14998.9995
$ws.Range("L134").Value = 22500
$ws.Range("M134").Value = -12463.9995
$ws.Range("N134").Value = -27570

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 5000
$ws.Range("I4").Value = 5000
$ws.Range("K4").Value = 5000
$ws.Range("M4").Value = -4888
# Row 7
$ws.Range("H7").Value = 86.92308
$ws.Range("I7").Value = 93.72727
$ws.Range("J7").Value = 49.5
$ws.Range("K7").Value = 93.72727
$ws.Range("L7").Value = 49.5
$ws.Range("M7").Value = 19.27273
$ws.Range("N7").Value = -275.5
# Row 31
$ws.Range("H31").Value = 4035.3333
$ws.Range("I31").Value = 2895
$ws.Range("K31").Value = 2895
$ws.Range("M31").Value = -2600
# Row 32
$ws.Range("H32").Value = 1530
$ws.Range("I32").Value = 1530
$ws.Range("K32").Value = 1530
$ws.Range("M32").Value = -1214
# Row 34
$ws.Range("H34").Value = 4035.3333
$ws.Range("I34").Value = 2895
$ws.Range("K34").Value = 2895
$ws.Range("M34").Value = -2693
# Row 99
$ws.Range("H99").Value = 13518
$ws.Range("I99").Value = 9411.5
$ws.Range("K99").Value = 9411.5
$ws.Range("M99").Value = -7913.5
# Row 105
$ws.Range("H105").Value = 2122.2222
$ws.Range("I105").Value = 575
$ws.Range("K105").Value = 575
$ws.Range("M105").Value = 1172
# Row 126
$ws.Range("H126").Value = 13518
$ws.Range("I126").Value = 9411.5
$ws.Range("K126").Value = 28234.5
$ws.Range("M126").Value = -25764.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 97
$ws.Range("H97").Value = 5532.125
$ws.Range("J97").Value = 5415.6
$ws.Range("L97").Value = 16246.8
$ws.Range("N97").Value = -17238.8
# Row 98
$ws.Range("H98").Value = 3762.4443
$ws.Range("J98").Value = 3512.5
$ws.Range("L98").Value = 10537.5
$ws.Range("N98").Value = -13533.5
# Row 115
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()
# Row 131
$ws.Range("H131").Value = 1777.0392
$ws.Range("I131").Value = 1570
$ws.Range("J131").Value = 1789.9791
$ws.Range("K131").Value = 4710
$ws.Range("L131").Value = 5369.9373
$ws.Range("M131").Value = 330
$ws.Range("N131").Value = -15449.9373
# Row 140
$ws.Range("H140").Value = 4467.75
$ws.Range("I140").Value = 4291.1665
$ws.Range("K140").Value = 12873.4995
$ws.Range("M140").Value = -7693.499500000002

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 129
$ws.Range("H129").Value = 45000
$ws.Range("J129").Value = 45000
$ws.Range("L129").Value = 45000
$ws.Range("N129").Value = -55000
# Row 130
$ws.Range("H130").Value = 85000
$ws.Range("J130").Value = 85000
$ws.Range("L130").Value = 85000
$ws.Range("N130").Value = -95040

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1075.4286
$ws.Range("I22").Value = 749
$ws.Range("J22").Value = 1129.8334
$ws.Range("K22").Value = 749
$ws.Range("L22").Value = 1129.8334
$ws.Range("M22").Value = -454
$ws.Range("N22").Value = -1719.8334
# Row 27
$ws.Range("H27").Value = 1075.4286
$ws.Range("I27").Value = 749
$ws.Range("J27").Value = 1129.8334
$ws.Range("K27").Value = 749
$ws.Range("L27").Value = 1129.8334
$ws.Range("M27").Value = -642
$ws.Range("N27").Value = -1343.8334
# Row 54
$ws.Range("H54").Value = 43666.332
$ws.Range("I54").Value = 43000
$ws.Range("J54").Value = 43999.5
$ws.Range("K54").Value = 43000
$ws.Range("L54").Value = 43999.5
$ws.Range("M54").Value = -42356
$ws.Range("N54").Value = -45287.5
# Row 68
$ws.Range("H68").Value = 4800.4443
$ws.Range("I68").Value = 4742.5713
$ws.Range("K68").Value = 4742.5713
$ws.Range("M68").Value = -3993.5713
# Row 71
$ws.Range("H71").Value = 4800.4443
$ws.Range("I71").Value = 4742.5713
$ws.Range("K71").Value = 23712.8565
$ws.Range("M71").Value = -19968.8565
# Row 82
$ws.Range("H82").Value = 1218.8889
$ws.Range("I82").Value = 661.1667
$ws.Range("K82").Value = 661.1667
$ws.Range("M82").Value = -300.1667
# Row 85
$ws.Range("H85").Value = 1218.8889
$ws.Range("I85").Value = 661.1667
$ws.Range("K85").Value = 661.1667
$ws.Range("M85").Value = 586.8333

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 525
$ws.Range("J2").Value = 525
$ws.Range("L2").Value = 525
$ws.Range("N2").Value = -749
# Row 4
$ws.Range("H4").Value = 6033.3335
$ws.Range("J4").Value = 537.5
$ws.Range("L4").Value = 537.5
$ws.Range("N4").Value = -763.5
# Row 96
$ws.Range("H96").Value = 4906.75
$ws.Range("I96").Value = 5186.6665
$ws.Range("K96").Value = 5186.6665
$ws.Range("M96").Value = -3813.6665

